$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990867465/1748730501814/koYLkezzyg-camera-video-segment-driver-1748730504314.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD77H2NTH%2F20250601%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250601T180021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBIaCXVzLXdlc3QtMiJHMEUCIAzQ7%2F3OWYJl0iHSY07s%2F1gMuzF4F4srHCirh4GxOP8TAiEA9%2BS9obOHOoT8Jbq68VYFFuM69DjNf5s3BfR6%2FTKn5Hcq5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBAuzViWLrC2GpXUrSq6A6hazXyRAMv72WYvwHgDQt3880vDJNQIIKHf1aSwITk5S1rJLlWMAxahW9XBDJINqIvComlxI0kXt04vK4K3YJSLLfflvWoRG2POR2nZhUDkSPEOBs6EpAK4WOEjDV93ezfa670AcnM%2Fappm5wHTAUh%2Bqr1WXBjNyu427yHMuRQV8sETJ3K5Q7zI51DYV%2BYtvo%2FkLSJB453C%2BbR7O95HbIFCDKkn8gJ0gB4wwbhN1ZCahm4ru91UZYK5gxhZm4%2Bsp9fcYPV8O%2FSPdVA0aSMmsVr%2FuIIz%2BsnhfUNWmAMicFD30yhT44D%2FWej1VPqLzDnyss2GOGzBARsn4iCtLLUwGVXond7BqOVhZE02EeEJrJCVA%2FSSqHMiT3W88zZ8kRqZntu4WcLXYr72I2cPALuqD7Db7kwE7sA4XhLSINJVepHfZ1OCL9ovGb0tW8%2FHXCYN6uBWQX9uHgPcVVwBuKsHOu46%2Bza1rcXxRIt1AJhKduu%2BZvaB79AGghsGoNeB5N6Sp3OWwNcuatdLmDx3wwum8h7Me9%2BdaAtGPzK%2Frgjfpz7D4qJbeoK4oMHfh11cqidznPk9kXL7cmRS8tcwuZ3ywQY6pQE3WPgkeXLon06KAHbX9TTgMD6dSjx3Nf%2Bry0VqYVhFJINj6eiaKKtPiffAVpqM6YoeQdEHwaKPTVSutnSvHt2b5HE4K501dL0GVq3bd1a93Eis%2FLCgvD623MBAtYstcYdidagk1Zq%2BJd%2B2TXPd74lLvgT12%2FyectW3tMOO2Qdvd3CFrrbKXeZmM3wyG%2FNwCk%2FodKVFh1MX01b2nHbvieH6EpfLmJs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2002%20Jun%202025%2002%3A00%3A21%20GMT&X-Amz-Signature=6d214129847b88f0d3fe16abad18d6411d39014e6cde6f1cb1b875e4b9846307'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990867465/1748730501814/oIElzMuidR-camera-video-segment-1748730504314.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD77H2NTH%2F20250601%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250601T180021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBIaCXVzLXdlc3QtMiJHMEUCIAzQ7%2F3OWYJl0iHSY07s%2F1gMuzF4F4srHCirh4GxOP8TAiEA9%2BS9obOHOoT8Jbq68VYFFuM69DjNf5s3BfR6%2FTKn5Hcq5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBAuzViWLrC2GpXUrSq6A6hazXyRAMv72WYvwHgDQt3880vDJNQIIKHf1aSwITk5S1rJLlWMAxahW9XBDJINqIvComlxI0kXt04vK4K3YJSLLfflvWoRG2POR2nZhUDkSPEOBs6EpAK4WOEjDV93ezfa670AcnM%2Fappm5wHTAUh%2Bqr1WXBjNyu427yHMuRQV8sETJ3K5Q7zI51DYV%2BYtvo%2FkLSJB453C%2BbR7O95HbIFCDKkn8gJ0gB4wwbhN1ZCahm4ru91UZYK5gxhZm4%2Bsp9fcYPV8O%2FSPdVA0aSMmsVr%2FuIIz%2BsnhfUNWmAMicFD30yhT44D%2FWej1VPqLzDnyss2GOGzBARsn4iCtLLUwGVXond7BqOVhZE02EeEJrJCVA%2FSSqHMiT3W88zZ8kRqZntu4WcLXYr72I2cPALuqD7Db7kwE7sA4XhLSINJVepHfZ1OCL9ovGb0tW8%2FHXCYN6uBWQX9uHgPcVVwBuKsHOu46%2Bza1rcXxRIt1AJhKduu%2BZvaB79AGghsGoNeB5N6Sp3OWwNcuatdLmDx3wwum8h7Me9%2BdaAtGPzK%2Frgjfpz7D4qJbeoK4oMHfh11cqidznPk9kXL7cmRS8tcwuZ3ywQY6pQE3WPgkeXLon06KAHbX9TTgMD6dSjx3Nf%2Bry0VqYVhFJINj6eiaKKtPiffAVpqM6YoeQdEHwaKPTVSutnSvHt2b5HE4K501dL0GVq3bd1a93Eis%2FLCgvD623MBAtYstcYdidagk1Zq%2BJd%2B2TXPd74lLvgT12%2FyectW3tMOO2Qdvd3CFrrbKXeZmM3wyG%2FNwCk%2FodKVFh1MX01b2nHbvieH6EpfLmJs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2002%20Jun%202025%2002%3A00%3A21%20GMT&X-Amz-Signature=365ead5cd70106f1ac6df3a0a65aac6b0c21e38ff7e016e4c21d3576350a0e4f'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990867465/1748728440718/4pjMy5ebbp-camera-video-segment-driver-1748728443218.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD77H2NTH%2F20250601%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250601T180021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBIaCXVzLXdlc3QtMiJHMEUCIAzQ7%2F3OWYJl0iHSY07s%2F1gMuzF4F4srHCirh4GxOP8TAiEA9%2BS9obOHOoT8Jbq68VYFFuM69DjNf5s3BfR6%2FTKn5Hcq5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBAuzViWLrC2GpXUrSq6A6hazXyRAMv72WYvwHgDQt3880vDJNQIIKHf1aSwITk5S1rJLlWMAxahW9XBDJINqIvComlxI0kXt04vK4K3YJSLLfflvWoRG2POR2nZhUDkSPEOBs6EpAK4WOEjDV93ezfa670AcnM%2Fappm5wHTAUh%2Bqr1WXBjNyu427yHMuRQV8sETJ3K5Q7zI51DYV%2BYtvo%2FkLSJB453C%2BbR7O95HbIFCDKkn8gJ0gB4wwbhN1ZCahm4ru91UZYK5gxhZm4%2Bsp9fcYPV8O%2FSPdVA0aSMmsVr%2FuIIz%2BsnhfUNWmAMicFD30yhT44D%2FWej1VPqLzDnyss2GOGzBARsn4iCtLLUwGVXond7BqOVhZE02EeEJrJCVA%2FSSqHMiT3W88zZ8kRqZntu4WcLXYr72I2cPALuqD7Db7kwE7sA4XhLSINJVepHfZ1OCL9ovGb0tW8%2FHXCYN6uBWQX9uHgPcVVwBuKsHOu46%2Bza1rcXxRIt1AJhKduu%2BZvaB79AGghsGoNeB5N6Sp3OWwNcuatdLmDx3wwum8h7Me9%2BdaAtGPzK%2Frgjfpz7D4qJbeoK4oMHfh11cqidznPk9kXL7cmRS8tcwuZ3ywQY6pQE3WPgkeXLon06KAHbX9TTgMD6dSjx3Nf%2Bry0VqYVhFJINj6eiaKKtPiffAVpqM6YoeQdEHwaKPTVSutnSvHt2b5HE4K501dL0GVq3bd1a93Eis%2FLCgvD623MBAtYstcYdidagk1Zq%2BJd%2B2TXPd74lLvgT12%2FyectW3tMOO2Qdvd3CFrrbKXeZmM3wyG%2FNwCk%2FodKVFh1MX01b2nHbvieH6EpfLmJs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2002%20Jun%202025%2002%3A00%3A21%20GMT&X-Amz-Signature=b1f0b24d6baf700cef30d5fcbb224aa702646e0f826768b79e18590d92d0e1f6'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990867465/1748728440718/jeBvGmmllL-camera-video-segment-1748728443218.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD77H2NTH%2F20250601%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250601T180021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBIaCXVzLXdlc3QtMiJHMEUCIAzQ7%2F3OWYJl0iHSY07s%2F1gMuzF4F4srHCirh4GxOP8TAiEA9%2BS9obOHOoT8Jbq68VYFFuM69DjNf5s3BfR6%2FTKn5Hcq5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBAuzViWLrC2GpXUrSq6A6hazXyRAMv72WYvwHgDQt3880vDJNQIIKHf1aSwITk5S1rJLlWMAxahW9XBDJINqIvComlxI0kXt04vK4K3YJSLLfflvWoRG2POR2nZhUDkSPEOBs6EpAK4WOEjDV93ezfa670AcnM%2Fappm5wHTAUh%2Bqr1WXBjNyu427yHMuRQV8sETJ3K5Q7zI51DYV%2BYtvo%2FkLSJB453C%2BbR7O95HbIFCDKkn8gJ0gB4wwbhN1ZCahm4ru91UZYK5gxhZm4%2Bsp9fcYPV8O%2FSPdVA0aSMmsVr%2FuIIz%2BsnhfUNWmAMicFD30yhT44D%2FWej1VPqLzDnyss2GOGzBARsn4iCtLLUwGVXond7BqOVhZE02EeEJrJCVA%2FSSqHMiT3W88zZ8kRqZntu4WcLXYr72I2cPALuqD7Db7kwE7sA4XhLSINJVepHfZ1OCL9ovGb0tW8%2FHXCYN6uBWQX9uHgPcVVwBuKsHOu46%2Bza1rcXxRIt1AJhKduu%2BZvaB79AGghsGoNeB5N6Sp3OWwNcuatdLmDx3wwum8h7Me9%2BdaAtGPzK%2Frgjfpz7D4qJbeoK4oMHfh11cqidznPk9kXL7cmRS8tcwuZ3ywQY6pQE3WPgkeXLon06KAHbX9TTgMD6dSjx3Nf%2Bry0VqYVhFJINj6eiaKKtPiffAVpqM6YoeQdEHwaKPTVSutnSvHt2b5HE4K501dL0GVq3bd1a93Eis%2FLCgvD623MBAtYstcYdidagk1Zq%2BJd%2B2TXPd74lLvgT12%2FyectW3tMOO2Qdvd3CFrrbKXeZmM3wyG%2FNwCk%2FodKVFh1MX01b2nHbvieH6EpfLmJs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2002%20Jun%202025%2002%3A00%3A21%20GMT&X-Amz-Signature=a6106ada723aa838f91996c7cf4ae992eb47cee8a7ba6b1c7a6d50ca883c08bf'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990867465/1748722729110/xV7hRByMSc-camera-video-segment-driver-1748722731610.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD77H2NTH%2F20250601%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250601T180021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBIaCXVzLXdlc3QtMiJHMEUCIAzQ7%2F3OWYJl0iHSY07s%2F1gMuzF4F4srHCirh4GxOP8TAiEA9%2BS9obOHOoT8Jbq68VYFFuM69DjNf5s3BfR6%2FTKn5Hcq5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBAuzViWLrC2GpXUrSq6A6hazXyRAMv72WYvwHgDQt3880vDJNQIIKHf1aSwITk5S1rJLlWMAxahW9XBDJINqIvComlxI0kXt04vK4K3YJSLLfflvWoRG2POR2nZhUDkSPEOBs6EpAK4WOEjDV93ezfa670AcnM%2Fappm5wHTAUh%2Bqr1WXBjNyu427yHMuRQV8sETJ3K5Q7zI51DYV%2BYtvo%2FkLSJB453C%2BbR7O95HbIFCDKkn8gJ0gB4wwbhN1ZCahm4ru91UZYK5gxhZm4%2Bsp9fcYPV8O%2FSPdVA0aSMmsVr%2FuIIz%2BsnhfUNWmAMicFD30yhT44D%2FWej1VPqLzDnyss2GOGzBARsn4iCtLLUwGVXond7BqOVhZE02EeEJrJCVA%2FSSqHMiT3W88zZ8kRqZntu4WcLXYr72I2cPALuqD7Db7kwE7sA4XhLSINJVepHfZ1OCL9ovGb0tW8%2FHXCYN6uBWQX9uHgPcVVwBuKsHOu46%2Bza1rcXxRIt1AJhKduu%2BZvaB79AGghsGoNeB5N6Sp3OWwNcuatdLmDx3wwum8h7Me9%2BdaAtGPzK%2Frgjfpz7D4qJbeoK4oMHfh11cqidznPk9kXL7cmRS8tcwuZ3ywQY6pQE3WPgkeXLon06KAHbX9TTgMD6dSjx3Nf%2Bry0VqYVhFJINj6eiaKKtPiffAVpqM6YoeQdEHwaKPTVSutnSvHt2b5HE4K501dL0GVq3bd1a93Eis%2FLCgvD623MBAtYstcYdidagk1Zq%2BJd%2B2TXPd74lLvgT12%2FyectW3tMOO2Qdvd3CFrrbKXeZmM3wyG%2FNwCk%2FodKVFh1MX01b2nHbvieH6EpfLmJs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2002%20Jun%202025%2002%3A00%3A21%20GMT&X-Amz-Signature=b60652e376025c6d83d5a1f7fdfc04490941987ce7d0e6277f752dd7cf006b2d'
$ws.Range("L4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990867465/1748722729110/Yc1IXI67aV-camera-video-segment-1748722731610.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD77H2NTH%2F20250601%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250601T180021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBIaCXVzLXdlc3QtMiJHMEUCIAzQ7%2F3OWYJl0iHSY07s%2F1gMuzF4F4srHCirh4GxOP8TAiEA9%2BS9obOHOoT8Jbq68VYFFuM69DjNf5s3BfR6%2FTKn5Hcq5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBAuzViWLrC2GpXUrSq6A6hazXyRAMv72WYvwHgDQt3880vDJNQIIKHf1aSwITk5S1rJLlWMAxahW9XBDJINqIvComlxI0kXt04vK4K3YJSLLfflvWoRG2POR2nZhUDkSPEOBs6EpAK4WOEjDV93ezfa670AcnM%2Fappm5wHTAUh%2Bqr1WXBjNyu427yHMuRQV8sETJ3K5Q7zI51DYV%2BYtvo%2FkLSJB453C%2BbR7O95HbIFCDKkn8gJ0gB4wwbhN1ZCahm4ru91UZYK5gxhZm4%2Bsp9fcYPV8O%2FSPdVA0aSMmsVr%2FuIIz%2BsnhfUNWmAMicFD30yhT44D%2FWej1VPqLzDnyss2GOGzBARsn4iCtLLUwGVXond7BqOVhZE02EeEJrJCVA%2FSSqHMiT3W88zZ8kRqZntu4WcLXYr72I2cPALuqD7Db7kwE7sA4XhLSINJVepHfZ1OCL9ovGb0tW8%2FHXCYN6uBWQX9uHgPcVVwBuKsHOu46%2Bza1rcXxRIt1AJhKduu%2BZvaB79AGghsGoNeB5N6Sp3OWwNcuatdLmDx3wwum8h7Me9%2BdaAtGPzK%2Frgjfpz7D4qJbeoK4oMHfh11cqidznPk9kXL7cmRS8tcwuZ3ywQY6pQE3WPgkeXLon06KAHbX9TTgMD6dSjx3Nf%2Bry0VqYVhFJINj6eiaKKtPiffAVpqM6YoeQdEHwaKPTVSutnSvHt2b5HE4K501dL0GVq3bd1a93Eis%2FLCgvD623MBAtYstcYdidagk1Zq%2BJd%2B2TXPd74lLvgT12%2FyectW3tMOO2Qdvd3CFrrbKXeZmM3wyG%2FNwCk%2FodKVFh1MX01b2nHbvieH6EpfLmJs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2002%20Jun%202025%2002%3A00%3A21%20GMT&X-Amz-Signature=d40dcd06af96194d0ca3635a68094d072af7dd39967551ca28e15566d5081155'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991154852/1748719417335/T9A5XCMwBf-camera-video-segment-driver-1748719419835.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD77H2NTH%2F20250601%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250601T180021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBIaCXVzLXdlc3QtMiJHMEUCIAzQ7%2F3OWYJl0iHSY07s%2F1gMuzF4F4srHCirh4GxOP8TAiEA9%2BS9obOHOoT8Jbq68VYFFuM69DjNf5s3BfR6%2FTKn5Hcq5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBAuzViWLrC2GpXUrSq6A6hazXyRAMv72WYvwHgDQt3880vDJNQIIKHf1aSwITk5S1rJLlWMAxahW9XBDJINqIvComlxI0kXt04vK4K3YJSLLfflvWoRG2POR2nZhUDkSPEOBs6EpAK4WOEjDV93ezfa670AcnM%2Fappm5wHTAUh%2Bqr1WXBjNyu427yHMuRQV8sETJ3K5Q7zI51DYV%2BYtvo%2FkLSJB453C%2BbR7O95HbIFCDKkn8gJ0gB4wwbhN1ZCahm4ru91UZYK5gxhZm4%2Bsp9fcYPV8O%2FSPdVA0aSMmsVr%2FuIIz%2BsnhfUNWmAMicFD30yhT44D%2FWej1VPqLzDnyss2GOGzBARsn4iCtLLUwGVXond7BqOVhZE02EeEJrJCVA%2FSSqHMiT3W88zZ8kRqZntu4WcLXYr72I2cPALuqD7Db7kwE7sA4XhLSINJVepHfZ1OCL9ovGb0tW8%2FHXCYN6uBWQX9uHgPcVVwBuKsHOu46%2Bza1rcXxRIt1AJhKduu%2BZvaB79AGghsGoNeB5N6Sp3OWwNcuatdLmDx3wwum8h7Me9%2BdaAtGPzK%2Frgjfpz7D4qJbeoK4oMHfh11cqidznPk9kXL7cmRS8tcwuZ3ywQY6pQE3WPgkeXLon06KAHbX9TTgMD6dSjx3Nf%2Bry0VqYVhFJINj6eiaKKtPiffAVpqM6YoeQdEHwaKPTVSutnSvHt2b5HE4K501dL0GVq3bd1a93Eis%2FLCgvD623MBAtYstcYdidagk1Zq%2BJd%2B2TXPd74lLvgT12%2FyectW3tMOO2Qdvd3CFrrbKXeZmM3wyG%2FNwCk%2FodKVFh1MX01b2nHbvieH6EpfLmJs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2002%20Jun%202025%2002%3A00%3A21%20GMT&X-Amz-Signature=32cd33118ad097c9b4ada0c4879fe3370c473d0954a88e44edb13c95d2a81e16'
$ws.Range("L5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991154852/1748719417335/iUK6XjjW8R-camera-video-segment-1748719419835.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD77H2NTH%2F20250601%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250601T180021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBIaCXVzLXdlc3QtMiJHMEUCIAzQ7%2F3OWYJl0iHSY07s%2F1gMuzF4F4srHCirh4GxOP8TAiEA9%2BS9obOHOoT8Jbq68VYFFuM69DjNf5s3BfR6%2FTKn5Hcq5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBAuzViWLrC2GpXUrSq6A6hazXyRAMv72WYvwHgDQt3880vDJNQIIKHf1aSwITk5S1rJLlWMAxahW9XBDJINqIvComlxI0kXt04vK4K3YJSLLfflvWoRG2POR2nZhUDkSPEOBs6EpAK4WOEjDV93ezfa670AcnM%2Fappm5wHTAUh%2Bqr1WXBjNyu427yHMuRQV8sETJ3K5Q7zI51DYV%2BYtvo%2FkLSJB453C%2BbR7O95HbIFCDKkn8gJ0gB4wwbhN1ZCahm4ru91UZYK5gxhZm4%2Bsp9fcYPV8O%2FSPdVA0aSMmsVr%2FuIIz%2BsnhfUNWmAMicFD30yhT44D%2FWej1VPqLzDnyss2GOGzBARsn4iCtLLUwGVXond7BqOVhZE02EeEJrJCVA%2FSSqHMiT3W88zZ8kRqZntu4WcLXYr72I2cPALuqD7Db7kwE7sA4XhLSINJVepHfZ1OCL9ovGb0tW8%2FHXCYN6uBWQX9uHgPcVVwBuKsHOu46%2Bza1rcXxRIt1AJhKduu%2BZvaB79AGghsGoNeB5N6Sp3OWwNcuatdLmDx3wwum8h7Me9%2BdaAtGPzK%2Frgjfpz7D4qJbeoK4oMHfh11cqidznPk9kXL7cmRS8tcwuZ3ywQY6pQE3WPgkeXLon06KAHbX9TTgMD6dSjx3Nf%2Bry0VqYVhFJINj6eiaKKtPiffAVpqM6YoeQdEHwaKPTVSutnSvHt2b5HE4K501dL0GVq3bd1a93Eis%2FLCgvD623MBAtYstcYdidagk1Zq%2BJd%2B2TXPd74lLvgT12%2FyectW3tMOO2Qdvd3CFrrbKXeZmM3wyG%2FNwCk%2FodKVFh1MX01b2nHbvieH6EpfLmJs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2002%20Jun%202025%2002%3A00%3A21%20GMT&X-Amz-Signature=022a8109534b5ed56a1d15df2f9b2e42e4849d1e88304d365ba4387eaa40d112'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991154852/1748708995160/wbBcub8c6J-camera-video-segment-driver-1748708997660.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD77H2NTH%2F20250601%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250601T180021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBIaCXVzLXdlc3QtMiJHMEUCIAzQ7%2F3OWYJl0iHSY07s%2F1gMuzF4F4srHCirh4GxOP8TAiEA9%2BS9obOHOoT8Jbq68VYFFuM69DjNf5s3BfR6%2FTKn5Hcq5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBAuzViWLrC2GpXUrSq6A6hazXyRAMv72WYvwHgDQt3880vDJNQIIKHf1aSwITk5S1rJLlWMAxahW9XBDJINqIvComlxI0kXt04vK4K3YJSLLfflvWoRG2POR2nZhUDkSPEOBs6EpAK4WOEjDV93ezfa670AcnM%2Fappm5wHTAUh%2Bqr1WXBjNyu427yHMuRQV8sETJ3K5Q7zI51DYV%2BYtvo%2FkLSJB453C%2BbR7O95HbIFCDKkn8gJ0gB4wwbhN1ZCahm4ru91UZYK5gxhZm4%2Bsp9fcYPV8O%2FSPdVA0aSMmsVr%2FuIIz%2BsnhfUNWmAMicFD30yhT44D%2FWej1VPqLzDnyss2GOGzBARsn4iCtLLUwGVXond7BqOVhZE02EeEJrJCVA%2FSSqHMiT3W88zZ8kRqZntu4WcLXYr72I2cPALuqD7Db7kwE7sA4XhLSINJVepHfZ1OCL9ovGb0tW8%2FHXCYN6uBWQX9uHgPcVVwBuKsHOu46%2Bza1rcXxRIt1AJhKduu%2BZvaB79AGghsGoNeB5N6Sp3OWwNcuatdLmDx3wwum8h7Me9%2BdaAtGPzK%2Frgjfpz7D4qJbeoK4oMHfh11cqidznPk9kXL7cmRS8tcwuZ3ywQY6pQE3WPgkeXLon06KAHbX9TTgMD6dSjx3Nf%2Bry0VqYVhFJINj6eiaKKtPiffAVpqM6YoeQdEHwaKPTVSutnSvHt2b5HE4K501dL0GVq3bd1a93Eis%2FLCgvD623MBAtYstcYdidagk1Zq%2BJd%2B2TXPd74lLvgT12%2FyectW3tMOO2Qdvd3CFrrbKXeZmM3wyG%2FNwCk%2FodKVFh1MX01b2nHbvieH6EpfLmJs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2002%20Jun%202025%2002%3A00%3A21%20GMT&X-Amz-Signature=ae232ca23abe698973c4586864653887fcf720b8adbb85aeed918f3b0bba314f'
$ws.Range("L6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991154852/1748708995160/JdTmyMZ5Lc-camera-video-segment-1748708997660.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD77H2NTH%2F20250601%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250601T180021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBIaCXVzLXdlc3QtMiJHMEUCIAzQ7%2F3OWYJl0iHSY07s%2F1gMuzF4F4srHCirh4GxOP8TAiEA9%2BS9obOHOoT8Jbq68VYFFuM69DjNf5s3BfR6%2FTKn5Hcq5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBAuzViWLrC2GpXUrSq6A6hazXyRAMv72WYvwHgDQt3880vDJNQIIKHf1aSwITk5S1rJLlWMAxahW9XBDJINqIvComlxI0kXt04vK4K3YJSLLfflvWoRG2POR2nZhUDkSPEOBs6EpAK4WOEjDV93ezfa670AcnM%2Fappm5wHTAUh%2Bqr1WXBjNyu427yHMuRQV8sETJ3K5Q7zI51DYV%2BYtvo%2FkLSJB453C%2BbR7O95HbIFCDKkn8gJ0gB4wwbhN1ZCahm4ru91UZYK5gxhZm4%2Bsp9fcYPV8O%2FSPdVA0aSMmsVr%2FuIIz%2BsnhfUNWmAMicFD30yhT44D%2FWej1VPqLzDnyss2GOGzBARsn4iCtLLUwGVXond7BqOVhZE02EeEJrJCVA%2FSSqHMiT3W88zZ8kRqZntu4WcLXYr72I2cPALuqD7Db7kwE7sA4XhLSINJVepHfZ1OCL9ovGb0tW8%2FHXCYN6uBWQX9uHgPcVVwBuKsHOu46%2Bza1rcXxRIt1AJhKduu%2BZvaB79AGghsGoNeB5N6Sp3OWwNcuatdLmDx3wwum8h7Me9%2BdaAtGPzK%2Frgjfpz7D4qJbeoK4oMHfh11cqidznPk9kXL7cmRS8tcwuZ3ywQY6pQE3WPgkeXLon06KAHbX9TTgMD6dSjx3Nf%2Bry0VqYVhFJINj6eiaKKtPiffAVpqM6YoeQdEHwaKPTVSutnSvHt2b5HE4K501dL0GVq3bd1a93Eis%2FLCgvD623MBAtYstcYdidagk1Zq%2BJd%2B2TXPd74lLvgT12%2FyectW3tMOO2Qdvd3CFrrbKXeZmM3wyG%2FNwCk%2FodKVFh1MX01b2nHbvieH6EpfLmJs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2002%20Jun%202025%2002%3A00%3A21%20GMT&X-Amz-Signature=b34d69c3c478696d694853eed2f52a82fa3f0171e08ca3769a64256805a8a7c0'
$ws.Range("K7").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474992807659/1748708610627/EAqUmKpcKh-camera-video-segment-driver-1748708613127.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD77H2NTH%2F20250601%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250601T180021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBIaCXVzLXdlc3QtMiJHMEUCIAzQ7%2F3OWYJl0iHSY07s%2F1gMuzF4F4srHCirh4GxOP8TAiEA9%2BS9obOHOoT8Jbq68VYFFuM69DjNf5s3BfR6%2FTKn5Hcq5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBAuzViWLrC2GpXUrSq6A6hazXyRAMv72WYvwHgDQt3880vDJNQIIKHf1aSwITk5S1rJLlWMAxahW9XBDJINqIvComlxI0kXt04vK4K3YJSLLfflvWoRG2POR2nZhUDkSPEOBs6EpAK4WOEjDV93ezfa670AcnM%2Fappm5wHTAUh%2Bqr1WXBjNyu427yHMuRQV8sETJ3K5Q7zI51DYV%2BYtvo%2FkLSJB453C%2BbR7O95HbIFCDKkn8gJ0gB4wwbhN1ZCahm4ru91UZYK5gxhZm4%2Bsp9fcYPV8O%2FSPdVA0aSMmsVr%2FuIIz%2BsnhfUNWmAMicFD30yhT44D%2FWej1VPqLzDnyss2GOGzBARsn4iCtLLUwGVXond7BqOVhZE02EeEJrJCVA%2FSSqHMiT3W88zZ8kRqZntu4WcLXYr72I2cPALuqD7Db7kwE7sA4XhLSINJVepHfZ1OCL9ovGb0tW8%2FHXCYN6uBWQX9uHgPcVVwBuKsHOu46%2Bza1rcXxRIt1AJhKduu%2BZvaB79AGghsGoNeB5N6Sp3OWwNcuatdLmDx3wwum8h7Me9%2BdaAtGPzK%2Frgjfpz7D4qJbeoK4oMHfh11cqidznPk9kXL7cmRS8tcwuZ3ywQY6pQE3WPgkeXLon06KAHbX9TTgMD6dSjx3Nf%2Bry0VqYVhFJINj6eiaKKtPiffAVpqM6YoeQdEHwaKPTVSutnSvHt2b5HE4K501dL0GVq3bd1a93Eis%2FLCgvD623MBAtYstcYdidagk1Zq%2BJd%2B2TXPd74lLvgT12%2FyectW3tMOO2Qdvd3CFrrbKXeZmM3wyG%2FNwCk%2FodKVFh1MX01b2nHbvieH6EpfLmJs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2002%20Jun%202025%2002%3A00%3A21%20GMT&X-Amz-Signature=187411f1c179c794e56b9d7fd172cb08b620c669e6893e36bed740d5243d08f4'
$ws.Range("L7").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474992807659/1748708610627/2sO4laNG63-camera-video-segment-1748708613127.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD77H2NTH%2F20250601%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250601T180021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBIaCXVzLXdlc3QtMiJHMEUCIAzQ7%2F3OWYJl0iHSY07s%2F1gMuzF4F4srHCirh4GxOP8TAiEA9%2BS9obOHOoT8Jbq68VYFFuM69DjNf5s3BfR6%2FTKn5Hcq5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBAuzViWLrC2GpXUrSq6A6hazXyRAMv72WYvwHgDQt3880vDJNQIIKHf1aSwITk5S1rJLlWMAxahW9XBDJINqIvComlxI0kXt04vK4K3YJSLLfflvWoRG2POR2nZhUDkSPEOBs6EpAK4WOEjDV93ezfa670AcnM%2Fappm5wHTAUh%2Bqr1WXBjNyu427yHMuRQV8sETJ3K5Q7zI51DYV%2BYtvo%2FkLSJB453C%2BbR7O95HbIFCDKkn8gJ0gB4wwbhN1ZCahm4ru91UZYK5gxhZm4%2Bsp9fcYPV8O%2FSPdVA0aSMmsVr%2FuIIz%2BsnhfUNWmAMicFD30yhT44D%2FWej1VPqLzDnyss2GOGzBARsn4iCtLLUwGVXond7BqOVhZE02EeEJrJCVA%2FSSqHMiT3W88zZ8kRqZntu4WcLXYr72I2cPALuqD7Db7kwE7sA4XhLSINJVepHfZ1OCL9ovGb0tW8%2FHXCYN6uBWQX9uHgPcVVwBuKsHOu46%2Bza1rcXxRIt1AJhKduu%2BZvaB79AGghsGoNeB5N6Sp3OWwNcuatdLmDx3wwum8h7Me9%2BdaAtGPzK%2Frgjfpz7D4qJbeoK4oMHfh11cqidznPk9kXL7cmRS8tcwuZ3ywQY6pQE3WPgkeXLon06KAHbX9TTgMD6dSjx3Nf%2Bry0VqYVhFJINj6eiaKKtPiffAVpqM6YoeQdEHwaKPTVSutnSvHt2b5HE4K501dL0GVq3bd1a93Eis%2FLCgvD623MBAtYstcYdidagk1Zq%2BJd%2B2TXPd74lLvgT12%2FyectW3tMOO2Qdvd3CFrrbKXeZmM3wyG%2FNwCk%2FodKVFh1MX01b2nHbvieH6EpfLmJs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2002%20Jun%202025%2002%3A00%3A21%20GMT&X-Amz-Signature=88a916b4f57382070fe528c66a9823dc89aa6607e58842cf40ed678f405d4126'
$ws.Range("K8").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474992807659/1748698902608/VYE0MXB4XQ-camera-video-segment-driver-1748698905108.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD77H2NTH%2F20250601%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250601T180021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBIaCXVzLXdlc3QtMiJHMEUCIAzQ7%2F3OWYJl0iHSY07s%2F1gMuzF4F4srHCirh4GxOP8TAiEA9%2BS9obOHOoT8Jbq68VYFFuM69DjNf5s3BfR6%2FTKn5Hcq5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBAuzViWLrC2GpXUrSq6A6hazXyRAMv72WYvwHgDQt3880vDJNQIIKHf1aSwITk5S1rJLlWMAxahW9XBDJINqIvComlxI0kXt04vK4K3YJSLLfflvWoRG2POR2nZhUDkSPEOBs6EpAK4WOEjDV93ezfa670AcnM%2Fappm5wHTAUh%2Bqr1WXBjNyu427yHMuRQV8sETJ3K5Q7zI51DYV%2BYtvo%2FkLSJB453C%2BbR7O95HbIFCDKkn8gJ0gB4wwbhN1ZCahm4ru91UZYK5gxhZm4%2Bsp9fcYPV8O%2FSPdVA0aSMmsVr%2FuIIz%2BsnhfUNWmAMicFD30yhT44D%2FWej1VPqLzDnyss2GOGzBARsn4iCtLLUwGVXond7BqOVhZE02EeEJrJCVA%2FSSqHMiT3W88zZ8kRqZntu4WcLXYr72I2cPALuqD7Db7kwE7sA4XhLSINJVepHfZ1OCL9ovGb0tW8%2FHXCYN6uBWQX9uHgPcVVwBuKsHOu46%2Bza1rcXxRIt1AJhKduu%2BZvaB79AGghsGoNeB5N6Sp3OWwNcuatdLmDx3wwum8h7Me9%2BdaAtGPzK%2Frgjfpz7D4qJbeoK4oMHfh11cqidznPk9kXL7cmRS8tcwuZ3ywQY6pQE3WPgkeXLon06KAHbX9TTgMD6dSjx3Nf%2Bry0VqYVhFJINj6eiaKKtPiffAVpqM6YoeQdEHwaKPTVSutnSvHt2b5HE4K501dL0GVq3bd1a93Eis%2FLCgvD623MBAtYstcYdidagk1Zq%2BJd%2B2TXPd74lLvgT12%2FyectW3tMOO2Qdvd3CFrrbKXeZmM3wyG%2FNwCk%2FodKVFh1MX01b2nHbvieH6EpfLmJs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2002%20Jun%202025%2002%3A00%3A21%20GMT&X-Amz-Signature=7f81e84ad3cff75a337823fcb9618af787a5d81af945a06e84c743b978fb870e'
$ws.Range("L8").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474992807659/1748698902608/AyCic9lIij-camera-video-segment-1748698905108.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD77H2NTH%2F20250601%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250601T180021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBIaCXVzLXdlc3QtMiJHMEUCIAzQ7%2F3OWYJl0iHSY07s%2F1gMuzF4F4srHCirh4GxOP8TAiEA9%2BS9obOHOoT8Jbq68VYFFuM69DjNf5s3BfR6%2FTKn5Hcq5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBAuzViWLrC2GpXUrSq6A6hazXyRAMv72WYvwHgDQt3880vDJNQIIKHf1aSwITk5S1rJLlWMAxahW9XBDJINqIvComlxI0kXt04vK4K3YJSLLfflvWoRG2POR2nZhUDkSPEOBs6EpAK4WOEjDV93ezfa670AcnM%2Fappm5wHTAUh%2Bqr1WXBjNyu427yHMuRQV8sETJ3K5Q7zI51DYV%2BYtvo%2FkLSJB453C%2BbR7O95HbIFCDKkn8gJ0gB4wwbhN1ZCahm4ru91UZYK5gxhZm4%2Bsp9fcYPV8O%2FSPdVA0aSMmsVr%2FuIIz%2BsnhfUNWmAMicFD30yhT44D%2FWej1VPqLzDnyss2GOGzBARsn4iCtLLUwGVXond7BqOVhZE02EeEJrJCVA%2FSSqHMiT3W88zZ8kRqZntu4WcLXYr72I2cPALuqD7Db7kwE7sA4XhLSINJVepHfZ1OCL9ovGb0tW8%2FHXCYN6uBWQX9uHgPcVVwBuKsHOu46%2Bza1rcXxRIt1AJhKduu%2BZvaB79AGghsGoNeB5N6Sp3OWwNcuatdLmDx3wwum8h7Me9%2BdaAtGPzK%2Frgjfpz7D4qJbeoK4oMHfh11cqidznPk9kXL7cmRS8tcwuZ3ywQY6pQE3WPgkeXLon06KAHbX9TTgMD6dSjx3Nf%2Bry0VqYVhFJINj6eiaKKtPiffAVpqM6YoeQdEHwaKPTVSutnSvHt2b5HE4K501dL0GVq3bd1a93Eis%2FLCgvD623MBAtYstcYdidagk1Zq%2BJd%2B2TXPd74lLvgT12%2FyectW3tMOO2Qdvd3CFrrbKXeZmM3wyG%2FNwCk%2FodKVFh1MX01b2nHbvieH6EpfLmJs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2002%20Jun%202025%2002%3A00%3A21%20GMT&X-Amz-Signature=dd81e0d04b57688c6dd1a83334b45d94cfd85f51e6c62a3648b5b54f7149cd9f'
$ws.Range("K9").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990867465/1748685998122/zLUrAcPiVy-camera-video-segment-driver-1748686000622.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD77H2NTH%2F20250601%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250601T180021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBIaCXVzLXdlc3QtMiJHMEUCIAzQ7%2F3OWYJl0iHSY07s%2F1gMuzF4F4srHCirh4GxOP8TAiEA9%2BS9obOHOoT8Jbq68VYFFuM69DjNf5s3BfR6%2FTKn5Hcq5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBAuzViWLrC2GpXUrSq6A6hazXyRAMv72WYvwHgDQt3880vDJNQIIKHf1aSwITk5S1rJLlWMAxahW9XBDJINqIvComlxI0kXt04vK4K3YJSLLfflvWoRG2POR2nZhUDkSPEOBs6EpAK4WOEjDV93ezfa670AcnM%2Fappm5wHTAUh%2Bqr1WXBjNyu427yHMuRQV8sETJ3K5Q7zI51DYV%2BYtvo%2FkLSJB453C%2BbR7O95HbIFCDKkn8gJ0gB4wwbhN1ZCahm4ru91UZYK5gxhZm4%2Bsp9fcYPV8O%2FSPdVA0aSMmsVr%2FuIIz%2BsnhfUNWmAMicFD30yhT44D%2FWej1VPqLzDnyss2GOGzBARsn4iCtLLUwGVXond7BqOVhZE02EeEJrJCVA%2FSSqHMiT3W88zZ8kRqZntu4WcLXYr72I2cPALuqD7Db7kwE7sA4XhLSINJVepHfZ1OCL9ovGb0tW8%2FHXCYN6uBWQX9uHgPcVVwBuKsHOu46%2Bza1rcXxRIt1AJhKduu%2BZvaB79AGghsGoNeB5N6Sp3OWwNcuatdLmDx3wwum8h7Me9%2BdaAtGPzK%2Frgjfpz7D4qJbeoK4oMHfh11cqidznPk9kXL7cmRS8tcwuZ3ywQY6pQE3WPgkeXLon06KAHbX9TTgMD6dSjx3Nf%2Bry0VqYVhFJINj6eiaKKtPiffAVpqM6YoeQdEHwaKPTVSutnSvHt2b5HE4K501dL0GVq3bd1a93Eis%2FLCgvD623MBAtYstcYdidagk1Zq%2BJd%2B2TXPd74lLvgT12%2FyectW3tMOO2Qdvd3CFrrbKXeZmM3wyG%2FNwCk%2FodKVFh1MX01b2nHbvieH6EpfLmJs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2002%20Jun%202025%2002%3A00%3A21%20GMT&X-Amz-Signature=d24da71645a64287f77cb5d80a83d4f91865636fde21bedb7a444203ca9985fb'
$ws.Range("L9").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990867465/1748685998122/r18KJ4zIte-camera-video-segment-1748686000622.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD77H2NTH%2F20250601%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250601T180021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBIaCXVzLXdlc3QtMiJHMEUCIAzQ7%2F3OWYJl0iHSY07s%2F1gMuzF4F4srHCirh4GxOP8TAiEA9%2BS9obOHOoT8Jbq68VYFFuM69DjNf5s3BfR6%2FTKn5Hcq5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBAuzViWLrC2GpXUrSq6A6hazXyRAMv72WYvwHgDQt3880vDJNQIIKHf1aSwITk5S1rJLlWMAxahW9XBDJINqIvComlxI0kXt04vK4K3YJSLLfflvWoRG2POR2nZhUDkSPEOBs6EpAK4WOEjDV93ezfa670AcnM%2Fappm5wHTAUh%2Bqr1WXBjNyu427yHMuRQV8sETJ3K5Q7zI51DYV%2BYtvo%2FkLSJB453C%2BbR7O95HbIFCDKkn8gJ0gB4wwbhN1ZCahm4ru91UZYK5gxhZm4%2Bsp9fcYPV8O%2FSPdVA0aSMmsVr%2FuIIz%2BsnhfUNWmAMicFD30yhT44D%2FWej1VPqLzDnyss2GOGzBARsn4iCtLLUwGVXond7BqOVhZE02EeEJrJCVA%2FSSqHMiT3W88zZ8kRqZntu4WcLXYr72I2cPALuqD7Db7kwE7sA4XhLSINJVepHfZ1OCL9ovGb0tW8%2FHXCYN6uBWQX9uHgPcVVwBuKsHOu46%2Bza1rcXxRIt1AJhKduu%2BZvaB79AGghsGoNeB5N6Sp3OWwNcuatdLmDx3wwum8h7Me9%2BdaAtGPzK%2Frgjfpz7D4qJbeoK4oMHfh11cqidznPk9kXL7cmRS8tcwuZ3ywQY6pQE3WPgkeXLon06KAHbX9TTgMD6dSjx3Nf%2Bry0VqYVhFJINj6eiaKKtPiffAVpqM6YoeQdEHwaKPTVSutnSvHt2b5HE4K501dL0GVq3bd1a93Eis%2FLCgvD623MBAtYstcYdidagk1Zq%2BJd%2B2TXPd74lLvgT12%2FyectW3tMOO2Qdvd3CFrrbKXeZmM3wyG%2FNwCk%2FodKVFh1MX01b2nHbvieH6EpfLmJs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2002%20Jun%202025%2002%3A00%3A21%20GMT&X-Amz-Signature=19d03bf806207296f2bb66d1b7a8a3a532db90a8a8a4bda6c93022268721703e'
$ws.Range("K10").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990867465/1748685979592/HA2n3X24AZ-camera-video-segment-driver-1748685982092.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD77H2NTH%2F20250601%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250601T180021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBIaCXVzLXdlc3QtMiJHMEUCIAzQ7%2F3OWYJl0iHSY07s%2F1gMuzF4F4srHCirh4GxOP8TAiEA9%2BS9obOHOoT8Jbq68VYFFuM69DjNf5s3BfR6%2FTKn5Hcq5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBAuzViWLrC2GpXUrSq6A6hazXyRAMv72WYvwHgDQt3880vDJNQIIKHf1aSwITk5S1rJLlWMAxahW9XBDJINqIvComlxI0kXt04vK4K3YJSLLfflvWoRG2POR2nZhUDkSPEOBs6EpAK4WOEjDV93ezfa670AcnM%2Fappm5wHTAUh%2Bqr1WXBjNyu427yHMuRQV8sETJ3K5Q7zI51DYV%2BYtvo%2FkLSJB453C%2BbR7O95HbIFCDKkn8gJ0gB4wwbhN1ZCahm4ru91UZYK5gxhZm4%2Bsp9fcYPV8O%2FSPdVA0aSMmsVr%2FuIIz%2BsnhfUNWmAMicFD30yhT44D%2FWej1VPqLzDnyss2GOGzBARsn4iCtLLUwGVXond7BqOVhZE02EeEJrJCVA%2FSSqHMiT3W88zZ8kRqZntu4WcLXYr72I2cPALuqD7Db7kwE7sA4XhLSINJVepHfZ1OCL9ovGb0tW8%2FHXCYN6uBWQX9uHgPcVVwBuKsHOu46%2Bza1rcXxRIt1AJhKduu%2BZvaB79AGghsGoNeB5N6Sp3OWwNcuatdLmDx3wwum8h7Me9%2BdaAtGPzK%2Frgjfpz7D4qJbeoK4oMHfh11cqidznPk9kXL7cmRS8tcwuZ3ywQY6pQE3WPgkeXLon06KAHbX9TTgMD6dSjx3Nf%2Bry0VqYVhFJINj6eiaKKtPiffAVpqM6YoeQdEHwaKPTVSutnSvHt2b5HE4K501dL0GVq3bd1a93Eis%2FLCgvD623MBAtYstcYdidagk1Zq%2BJd%2B2TXPd74lLvgT12%2FyectW3tMOO2Qdvd3CFrrbKXeZmM3wyG%2FNwCk%2FodKVFh1MX01b2nHbvieH6EpfLmJs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2002%20Jun%202025%2002%3A00%3A21%20GMT&X-Amz-Signature=03a0ad23d26aa7fb200e7e21be7a904ab81d660cd0ab84bc4054da8addde1414'
$ws.Range("L10").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990867465/1748685979592/icdf7xm5Gw-camera-video-segment-1748685982092.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD77H2NTH%2F20250601%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250601T180021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBIaCXVzLXdlc3QtMiJHMEUCIAzQ7%2F3OWYJl0iHSY07s%2F1gMuzF4F4srHCirh4GxOP8TAiEA9%2BS9obOHOoT8Jbq68VYFFuM69DjNf5s3BfR6%2FTKn5Hcq5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBAuzViWLrC2GpXUrSq6A6hazXyRAMv72WYvwHgDQt3880vDJNQIIKHf1aSwITk5S1rJLlWMAxahW9XBDJINqIvComlxI0kXt04vK4K3YJSLLfflvWoRG2POR2nZhUDkSPEOBs6EpAK4WOEjDV93ezfa670AcnM%2Fappm5wHTAUh%2Bqr1WXBjNyu427yHMuRQV8sETJ3K5Q7zI51DYV%2BYtvo%2FkLSJB453C%2BbR7O95HbIFCDKkn8gJ0gB4wwbhN1ZCahm4ru91UZYK5gxhZm4%2Bsp9fcYPV8O%2FSPdVA0aSMmsVr%2FuIIz%2BsnhfUNWmAMicFD30yhT44D%2FWej1VPqLzDnyss2GOGzBARsn4iCtLLUwGVXond7BqOVhZE02EeEJrJCVA%2FSSqHMiT3W88zZ8kRqZntu4WcLXYr72I2cPALuqD7Db7kwE7sA4XhLSINJVepHfZ1OCL9ovGb0tW8%2FHXCYN6uBWQX9uHgPcVVwBuKsHOu46%2Bza1rcXxRIt1AJhKduu%2BZvaB79AGghsGoNeB5N6Sp3OWwNcuatdLmDx3wwum8h7Me9%2BdaAtGPzK%2Frgjfpz7D4qJbeoK4oMHfh11cqidznPk9kXL7cmRS8tcwuZ3ywQY6pQE3WPgkeXLon06KAHbX9TTgMD6dSjx3Nf%2Bry0VqYVhFJINj6eiaKKtPiffAVpqM6YoeQdEHwaKPTVSutnSvHt2b5HE4K501dL0GVq3bd1a93Eis%2FLCgvD623MBAtYstcYdidagk1Zq%2BJd%2B2TXPd74lLvgT12%2FyectW3tMOO2Qdvd3CFrrbKXeZmM3wyG%2FNwCk%2FodKVFh1MX01b2nHbvieH6EpfLmJs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2002%20Jun%202025%2002%3A00%3A21%20GMT&X-Amz-Signature=681f0a3cd4ab37d06324a02e03cbaf29b296fde505421ae52ce1284fb5b3636c'
$ws.Range("K11").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990867465/1748684794984/nNtvjgl14a-camera-video-segment-driver-1748684797484.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD77H2NTH%2F20250601%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250601T180021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBIaCXVzLXdlc3QtMiJHMEUCIAzQ7%2F3OWYJl0iHSY07s%2F1gMuzF4F4srHCirh4GxOP8TAiEA9%2BS9obOHOoT8Jbq68VYFFuM69DjNf5s3BfR6%2FTKn5Hcq5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBAuzViWLrC2GpXUrSq6A6hazXyRAMv72WYvwHgDQt3880vDJNQIIKHf1aSwITk5S1rJLlWMAxahW9XBDJINqIvComlxI0kXt04vK4K3YJSLLfflvWoRG2POR2nZhUDkSPEOBs6EpAK4WOEjDV93ezfa670AcnM%2Fappm5wHTAUh%2Bqr1WXBjNyu427yHMuRQV8sETJ3K5Q7zI51DYV%2BYtvo%2FkLSJB453C%2BbR7O95HbIFCDKkn8gJ0gB4wwbhN1ZCahm4ru91UZYK5gxhZm4%2Bsp9fcYPV8O%2FSPdVA0aSMmsVr%2FuIIz%2BsnhfUNWmAMicFD30yhT44D%2FWej1VPqLzDnyss2GOGzBARsn4iCtLLUwGVXond7BqOVhZE02EeEJrJCVA%2FSSqHMiT3W88zZ8kRqZntu4WcLXYr72I2cPALuqD7Db7kwE7sA4XhLSINJVepHfZ1OCL9ovGb0tW8%2FHXCYN6uBWQX9uHgPcVVwBuKsHOu46%2Bza1rcXxRIt1AJhKduu%2BZvaB79AGghsGoNeB5N6Sp3OWwNcuatdLmDx3wwum8h7Me9%2BdaAtGPzK%2Frgjfpz7D4qJbeoK4oMHfh11cqidznPk9kXL7cmRS8tcwuZ3ywQY6pQE3WPgkeXLon06KAHbX9TTgMD6dSjx3Nf%2Bry0VqYVhFJINj6eiaKKtPiffAVpqM6YoeQdEHwaKPTVSutnSvHt2b5HE4K501dL0GVq3bd1a93Eis%2FLCgvD623MBAtYstcYdidagk1Zq%2BJd%2B2TXPd74lLvgT12%2FyectW3tMOO2Qdvd3CFrrbKXeZmM3wyG%2FNwCk%2FodKVFh1MX01b2nHbvieH6EpfLmJs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2002%20Jun%202025%2002%3A00%3A21%20GMT&X-Amz-Signature=e08d2b788c46785c24d1f62019286544adfe08a4d5f8e95603313e7978074003'
$ws.Range("L11").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990867465/1748684794984/cH6NSyWnjp-camera-video-segment-1748684797484.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSD77H2NTH%2F20250601%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250601T180021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBIaCXVzLXdlc3QtMiJHMEUCIAzQ7%2F3OWYJl0iHSY07s%2F1gMuzF4F4srHCirh4GxOP8TAiEA9%2BS9obOHOoT8Jbq68VYFFuM69DjNf5s3BfR6%2FTKn5Hcq5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBAuzViWLrC2GpXUrSq6A6hazXyRAMv72WYvwHgDQt3880vDJNQIIKHf1aSwITk5S1rJLlWMAxahW9XBDJINqIvComlxI0kXt04vK4K3YJSLLfflvWoRG2POR2nZhUDkSPEOBs6EpAK4WOEjDV93ezfa670AcnM%2Fappm5wHTAUh%2Bqr1WXBjNyu427yHMuRQV8sETJ3K5Q7zI51DYV%2BYtvo%2FkLSJB453C%2BbR7O95HbIFCDKkn8gJ0gB4wwbhN1ZCahm4ru91UZYK5gxhZm4%2Bsp9fcYPV8O%2FSPdVA0aSMmsVr%2FuIIz%2BsnhfUNWmAMicFD30yhT44D%2FWej1VPqLzDnyss2GOGzBARsn4iCtLLUwGVXond7BqOVhZE02EeEJrJCVA%2FSSqHMiT3W88zZ8kRqZntu4WcLXYr72I2cPALuqD7Db7kwE7sA4XhLSINJVepHfZ1OCL9ovGb0tW8%2FHXCYN6uBWQX9uHgPcVVwBuKsHOu46%2Bza1rcXxRIt1AJhKduu%2BZvaB79AGghsGoNeB5N6Sp3OWwNcuatdLmDx3wwum8h7Me9%2BdaAtGPzK%2Frgjfpz7D4qJbeoK4oMHfh11cqidznPk9kXL7cmRS8tcwuZ3ywQY6pQE3WPgkeXLon06KAHbX9TTgMD6dSjx3Nf%2Bry0VqYVhFJINj6eiaKKtPiffAVpqM6YoeQdEHwaKPTVSutnSvHt2b5HE4K501dL0GVq3bd1a93Eis%2FLCgvD623MBAtYstcYdidagk1Zq%2BJd%2B2TXPd74lLvgT12%2FyectW3tMOO2Qdvd3CFrrbKXeZmM3wyG%2FNwCk%2FodKVFh1MX01b2nHbvieH6EpfLmJs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2002%20Jun%202025%2002%3A00%3A21%20GMT&X-Amz-Signature=f2aee329399dcea950f1606a58cce6ff2d9df931ffb40f0a2abde5aa1887b9eb'
